$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "I3"   = 2.5
    "H6"   = 2.1
    "P6"   = 1.62
    "F7"   = 2.24
    "G7"   = 2.7
    "H7"   = 2.96
    "I7"   = 3.45
    "J7"   = 2.94
    "F8"   = 4.4
    "I8"   = 1.65
    "J8"   = 4.3
    "F9"   = 1.92
    "G9"   = 2.16
    "H9"   = 3.45
    "I9"   = 4.5
    "J9"   = 3.8
    "K9"   = 4.8
    "Q9"   = 1.58
    "G10"  = 5.5
    "H10"  = 1.68
    "I10"  = 1.7
    "N10"  = 4.9
    "P10"  = 2.34
    "S10"  = 2.7
    "T10"  = 1.73
    "U10"  = 2.18
    "V10"  = 2.42
    "W10"  = 1.22
    "AO10" = 8.4
    "L11"  = 1.42
    "P11"  = 1.93
    "Y11"  = 17.5
    "AJ11" = 38
    "AN11" = 19
    "Q13"  = 1.78
    "J14"  = 3.1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
